# "Exit info list.xlsx" update:
#  - adds a row for the new "Transfer" motor, a new "zAvit" pot/angle row,
#    and fills in the B/C/D (Device/Function/Slot) columns for several
#    existing rows that previously only had a Name in column A.
#  - row 19 (old "Cannon" micro-switch entry) is vacated; its data (plus the
#    other micro-switch rows) is rewritten one row lower (20-25) so a blank
#    spacer row still separates the "Jaguar Motor" block from the
#    "Micro Switch" block.
#  - widens column C to fit the new, longer "Function" descriptions.
#  - leaves the active selection on C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- existing header / drive / joystick rows (text unchanged, kept for clarity) ---
$ws.Range("A2").Value = "refnum name"
$ws.Range("B2").Value = "Device"
$ws.Range("C2").Value = "Function "
$ws.Range("D2").Value = "Slot"

$ws.Range("A4").Value = "Drive"
$ws.Range("B4").Value = "4 Motors"
$ws.Range("C4").Value = "drive"
$ws.Range("D4").Value = "PWM 1-4"

$ws.Range("A5").Value = "Pilot"
$ws.Range("B5").Value = "Joystick"
$ws.Range("C5").Value = "driver's joystick"
$ws.Range("D5").Value = "USB 1"

$ws.Range("A6").Value = "CoPilot"
$ws.Range("B6").Value = "Joystick"
$ws.Range("C6").Value = "operator's joystick"
$ws.Range("D6").Value = "USB 2"

# --- cannon / shooter motors ---
$ws.Range("A8").Value = "First Cannon"
$ws.Range("B8").Value = "Jaguar Motor"
$ws.Range("C8").Value = "outer shooting cannon"
$ws.Range("D8").Value = "PWM 7"

$ws.Range("A9").Value = "Second Cannon"
$ws.Range("B9").Value = "Jaguar Motor"
$ws.Range("C9").Value = "inner shooting cannon"
$ws.Range("D9").Value = "PWM 8"

$ws.Range("A10").Value = "zAvit"
$ws.Range("B10").Value = "Jaguar Motor"
$ws.Range("C10").Value = "cannon angle"
$ws.Range("D10").Value = "PWM 9"

# --- sensors ---
$ws.Range("A12").Value = "CannonAI"
$ws.Range("B12").Value = "Potentiometer"
$ws.Range("C12").Value = "measures angle of cannon"
$ws.Range("D12").Value = "AI 1"

$ws.Range("A13").Value = "Enc"
$ws.Range("B13").Value = "Encoder"
$ws.Range("C13").Value = "measures shooting motor's RPM"
$ws.Range("D13").Value = "DIO 7-8"

# --- frisbee-handling motors ---
$ws.Range("A15").Value = "Flipper"
$ws.Range("B15").Value = "Jaguar Motor"
$ws.Range("C15").Value = "flipper's motor"
$ws.Range("D15").Value = "PWM 9"

$ws.Range("A16").Value = "Conveyer Belt"
$ws.Range("B16").Value = "Jaguar Motor"
$ws.Range("D16").Value = "PWM 6"

$ws.Range("A17").Value = "Roller"
$ws.Range("B17").Value = "Jaguar Motor"
$ws.Range("D17").Value = "PWM 5"

$ws.Range("A18").Value = "Transfer"
$ws.Range("B18").Value = "Jaguar Motor"
$ws.Range("D18").Value = "PWM 10"

# row 19 is now the blank spacer row before the micro-switch block - make
# sure the old "Cannon" value that used to live here is gone.
$ws.Range("A19").ClearContents()

# --- micro switches (moved down one row, 19->20 ... 24->25) ---
$ws.Range("A20").Value = "Cannon"
$ws.Range("B20").Value = "Micro Switch"
$ws.Range("D20").Value = "DIO 5"

$ws.Range("A21").Value = "InFlipper"
$ws.Range("B21").Value = "Micro Switch"
$ws.Range("D21").Value = "DIO 4"

$ws.Range("A22").Value = "FlipperAtPlace"
$ws.Range("B22").Value = "Micro Switch"
$ws.Range("D22").Value = "DIO 3"

$ws.Range("A23").Value = "Middle"
$ws.Range("B23").Value = "Micro Switch"
$ws.Range("D23").Value = "DIO 2"

$ws.Range("A24").Value = "Bottom"
$ws.Range("B24").Value = "Micro Switch"
$ws.Range("D24").Value = "DIO 1"

$ws.Range("A25").Value = "FrisbeeDirection"
$ws.Range("B25").Value = "Micro Switch"
$ws.Range("D25").Value = "DIO 6"

# --- column C needs to be noticeably wider for the new descriptions ---
$ws.Columns.Item(3).ColumnWidth = 26.85

# --- leave the selection where the author left it ---
$ws.Range("C17").Select()
